# Generate Report for Handoff
#
# A new source file (c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md) has been
# picked up by the localization pipeline. It is inserted as a new row,
# directly above the existing ".localization-config" bookkeeping row, on
# every report sheet (Overview, zh-cn, de-de). The old
# ".localization-config" row is pushed down by one row.

$wb = $excel.ActiveWorkbook

$newFile        = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.md"
$repoCommit     = "c463897b4bc5670b746f73360af9d2c57a535b03"
$newFileUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/e2e/$newFile"
$configUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/.localization-config"

$zhXlf          = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.zh-cn.xlf"
$zhXlfCommit    = "8341175a15e72d81211d1263b4bae898989241a1"
$zhXlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhXlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$zhDatetime     = "2016-03-10 03:10:20"

$deXlf          = "c2ffc5e8-7a8c-401d-b5a6-c6bc549c330b.a4fb2c484b80696ce0799656c3a2fb538b02d0f1.de-de.xlf"
$deXlfCommit    = "d94394ebae535fe0b9027c2b4b735700a647dc7a"
$deXlfUrl       = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deXlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"
$deDatetime     = "2016-03-10 03:10:23"

$readyStatus    = "Ready for handoff"
$configStatus   = "Not to be localized"
$epoch          = "0001-01-01 00:00:00"
$includeReason  = "Include"
$ignoredReason  = "Ignored"

# ----------------------------------------------------------------------
# Sheet "Overview" (3 columns: File Name | zh-cn | de-de)
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop whatever hyperlinks exist so stale ones don't linger on the wrong
# row once the old ".localization-config" row is shifted down.
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus

$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = $configStatus
$wsOverview.Range("C4").Value = $configStatus

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newFileUrl.Replace($newFile, "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md"), "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", ".localization-config") | Out-Null

# ----------------------------------------------------------------------
# Sheets "zh-cn" / "de-de" (9 columns, see table header row)
# ----------------------------------------------------------------------
function Fill-LangSheet($ws, $origMdUrl, $origXlfUrl, $origXlfDisplay, $origDatetime, $newXlfUrl, $newXlfDisplay, $newDatetime) {
    $ws.Range("A1").Hyperlinks.Delete()

    # Row 2 keeps the original source file, unchanged.
    $ws.Range("A2").Value = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md"
    $ws.Range("B2").Value = $readyStatus
    $ws.Range("C2").Value = $origXlfDisplay
    $ws.Range("D2").Value = $origDatetime
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = $includeReason

    # Row 3: the newly picked-up source file.
    $ws.Range("A3").Value = $newFile
    $ws.Range("B3").Value = $readyStatus
    $ws.Range("C3").Value = $newXlfDisplay
    $ws.Range("D3").Value = $newDatetime
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = $includeReason

    # Row 4: the ".localization-config" bookkeeping row, pushed down.
    $ws.Range("A4").Value = ".localization-config"
    $ws.Range("B4").Value = $configStatus
    $ws.Range("D4").Value = $epoch
    $ws.Range("G4").Value = $epoch
    $ws.Range("H4").Value = $ignoredReason

    $ws.Hyperlinks.Add($ws.Range("A2"), $origMdUrl, "", "", "17d16921-d803-4efd-bb22-ef1ca06a2a3e.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $origXlfUrl, "", "", $origXlfDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $newFileUrl, "", "", $newFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $newXlfUrl, "", "", $newXlfDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", ".localization-config") | Out-Null
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$origZhXlf = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.zh-cn.xlf"
$origZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8341175a15e72d81211d1263b4bae898989241a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$origZhXlf"
Fill-LangSheet $wsZh "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md" $origZhXlfUrl $origZhXlf "2016-03-10 03:09:36" $zhXlfUrl $zhXlf $zhDatetime

$wsDe = $wb.Worksheets.Item("de-de")
$origDeXlf = "17d16921-d803-4efd-bb22-ef1ca06a2a3e.fe2e2d909773b74bf09afd118f409a254cc6e502.de-de.xlf"
$origDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d94394ebae535fe0b9027c2b4b735700a647dc7a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$origDeXlf"
Fill-LangSheet $wsDe "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/e2e/17d16921-d803-4efd-bb22-ef1ca06a2a3e.md" $origDeXlfUrl $origDeXlf "2016-03-10 03:09:39" $deXlfUrl $deXlf $deDatetime
